# Gantt chart workbook: add the new "Designed a Main Page for the website
# to use(Craig)" task row under the existing task list on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task name (Column A)
$ws.Range("A11").Value = "Designed a Main Page for the website to use(Craig)"

# Start date for the new task (Column B) - 1/24/2021, matching the date
# format already used by the rows above it.
$ws.Range("B11").Value = 44220
$ws.Range("B11").NumberFormat = "d-mmm"

# Days taken (Column C)
$ws.Range("C11").Value = 1

# Leave the selection on the newly added row, like the author did.
$ws.Range("A11").Select()
